# Feedback fixes batch 2
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$subjectText = "Feedback fix batch 2"
$descriptionText = "I started with creating a new service, the tagService, and putting the tag related stuff from facility service in there. This was very simple to do. Then i fixed the issue where, when searching with a filter for tags, the API would only return the tag you were looking for, not the other tags that were part of the facilities you fetched. Its also being returned as an array, so i think the tags in a string issue is returned everywhere now."

# Row 24: Subject / Hours / Date / Description
$ws.Range("A24").Value = $subjectText
$ws.Range("B24").Value = 1
$ws.Range("C24").Value = [DateTime]::FromOADate(45796)
$ws.Range("D24").Value = $descriptionText

# Match the style of the row above (D column uses wrap text)
$ws.Range("D24").WrapText = $true
$ws.Rows.Item(24).RowHeight = 39.75

# Update the current selection like in the diff
$ws.Range("D26").Select()
